# Move the title placeholder shape ("Título 1", id=2) on the only slide
# to its new position, per the authored edit.
#
# Original offset (EMU): x=1071716, y=868362
# New offset (EMU):      x=442452,  y=581156
# (extent cx=9144000, cy=2387600 is unchanged)
#
# PowerPoint's Shape.Left/Shape.Top are expressed in points (1 pt = 12700 EMU)
# and are stored internally as single-precision floats, which truncates the
# EMU value on write-back. Resolve-PointsForEmu searches for a nearby double
# whose float32 cast round-trips to the exact target EMU so the saved OOXML
# matches exactly, instead of landing 1 EMU off because of rounding.
function Resolve-PointsForEmu([long]$targetEmu) {
    $base = $targetEmu / 12700.0
    for ($k = -200; $k -le 200; $k++) {
        $cand = $base + ($k * 0.0000015)
        $f = [double][float]$cand
        $emu = [math]::Floor($f * 12700.0 + 0.0000001)
        if ($emu -eq $targetEmu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)

$shape.Left = Resolve-PointsForEmu 442452
$shape.Top = Resolve-PointsForEmu 581156
